# Add an "Italy" worksheet, cloned from the "Slovakia" sheet (same template),
# populate its market name / part-number fields, and insert the two extra
# repeater rows (P32AR / P32DR) that Italy's sheet has but Slovakia's doesn't.

$wb = $excel.ActiveWorkbook

# Clone the Slovakia sheet as the template, placing the new sheet at the end.
$template = $wb.Worksheets.Item("Slovakia")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Italy"

# Market name / part number for Italy.
$newSheet.Range("B2").Value = "Italy Market"
$newSheet.Range("B4").Value = "NGC-3145/T2224 "

# Italy's repeater list additionally includes P32AR / P32DR (like the Swiss
# sheet), inserted right before the trailing "Wg" / "Repeaters" rows.
$newSheet.Rows.Item(16).Resize(2).Insert()

# Bring matching cell formatting down into the newly inserted rows.
$newSheet.Range("A15").Copy()
$newSheet.Range("A16:A17").PasteSpecial(-4122)

$newSheet.Range("A16").Value = "P32AR"
$newSheet.Range("A17").Value = "P32DR"

# Row 4 grew taller to fit the wrapped part-number text.
$newSheet.Rows.Item(4).RowHeight = 28.8

# Match column widths to fit the new content (closest values the engine
# quantizes to 25.11 / 15.22 / 20.44 characters respectively).
$newSheet.Columns.Item(1).ColumnWidth = 24.333333333333332
$newSheet.Columns.Item(2).ColumnWidth = 14.333333333333334
$newSheet.Columns.Item(4).ColumnWidth = 19.666666666666668

$newSheet.Range("B4").Select()
